$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'69.158.85"
$ws.Range("E2").Value = "  -1.32%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.492.58"
$ws.Range("E3").Value = "  -2.79%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'573.40"
$ws.Range("E5").Value = "  -1.27%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'185.04"
$ws.Range("E6").Value = "  -3.25%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "'3.482.91"
$ws.Range("E7").Value = "  -2.99%  "

# Row 8 - XRP
$ws.Range("D8").Value = "'0.610"
$ws.Range("E8").Value = "  -3.58%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.01%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +2.61%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -2.89%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "'53.97"
$ws.Range("E12").Value = "  -3.50%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -2.36%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -3.02%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "'4.050.43"
$ws.Range("E15").Value = "  -2.99%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "'19.27"
$ws.Range("E16").Value = "  -3.63%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "'69.091.52"
$ws.Range("E17").Value = "  -1.44%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "'3.495.55"
$ws.Range("E18").Value = "  -2.79%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "'12.23"
$ws.Range("E19").Value = "  -3.78%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  -1.21%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'540.70"
$ws.Range("E21").Value = "  +12.28%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  -3.88%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("D23").Value = "'18.36"
$ws.Range("E23").Value = "  -4.61%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "'4.94"
$ws.Range("E24").Value = "  -1.89%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "'4.41"
$ws.Range("E25").Value = "  +0.00%  "

# Row 26 - Litecoin
$ws.Range("D26").Value = "'93.39"
$ws.Range("E26").Value = "  -2.14%  "

# Row 27 - RenderToken
$ws.Range("D27").Value = "'11.26"
$ws.Range("E27").Value = "  +1.16%  "

# Row 28 - ImmutableX
$ws.Range("E28").Value = "  -2.10%  "

# Row 29 - Filecoin
$ws.Range("D29").Value = "'9.07"
$ws.Range("E29").Value = "  -3.81%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "'31.70"
$ws.Range("E30").Value = "  -1.50%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "'7.23"
$ws.Range("E31").Value = "  -5.73%  "

# Row 32 - Cosmos
$ws.Range("D32").Value = "'12.57"
$ws.Range("E32").Value = "  +2.76%  "

# Row 33 - OKB
$ws.Range("D33").Value = "'64.29"
$ws.Range("E33").Value = "  -3.68%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  -5.67%  "

# Row 35 - Bittensor
$ws.Range("D35").Value = "'532.61"
$ws.Range("E35").Value = "  -8.78%  "

# Row 36 - Fetch.AI
$ws.Range("D36").Value = "'3.08"
$ws.Range("E36").Value = "  +8.36%  "

# Row 37 - InjectiveProtocol
$ws.Range("D37").Value = "'37.80"
$ws.Range("E37").Value = "  -3.11%  "

# Row 38 - TheGraph/Dai swap
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.09%  "

# Row 39 - Dai/TheGraph swap
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "'0.399"
$ws.Range("E39").Value = "  +0.39%  "

# Row 40 - PEPE
$ws.Range("D40").Value = "'0.0₃0759"
$ws.Range("E40").Value = "  -5.66%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  -3.05%  "

# Row 42 - Kaspa
$ws.Range("D42").Value = "'0.132"
$ws.Range("E42").Value = "  -2.82%  "

# Row 43 - Maker
$ws.Range("D43").Value = "'3.290.68"
$ws.Range("E43").Value = "  +2.02%  "

# Row 44 - dogwifhat
$ws.Range("D44").Value = "'3.02"
$ws.Range("E44").Value = "  -9.02%  "

# Row 45 - ThetaToken
$ws.Range("E45").Value = "  -3.43%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  -1.72%  "

# Row 47 - ApeXProtocol
$ws.Range("E47").Value = "  +4.47%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  -3.56%  "

# Row 49 - THORChain
$ws.Range("D49").Value = "'8.88"
$ws.Range("E49").Value = "  -6.48%  "

# Row 50 - FirstDigitalUSD
$ws.Range("D50").Value = "'0.998"
$ws.Range("E50").Value = "  -0.10%  "

# Row 51 - Monero
$ws.Range("D51").Value = "'136.79"
$ws.Range("E51").Value = "  +1.84%  "
